$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.145.44"
$ws.Range("E2").Value = "  +2.83%  "

$ws.Range("D3").Value = "2.281.27"
$ws.Range("E3").Value = "  +2.89%  "

$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").Value = "319.14"
$ws.Range("E5").Value = "  +1.23%  "

$ws.Range("D6").Value = "107.31"
$ws.Range("E6").Value = "  +8.40%  "

$ws.Range("E7").Value = "  +0.51%  "

$ws.Range("D9").Value = "0.575"
$ws.Range("E9").Value = "  +2.55%  "

$ws.Range("D10").Value = "39.17"
$ws.Range("E10").Value = "  +7.06%  "

$ws.Range("E11").Value = "  +2.04%  "

$ws.Range("D12").Value = "7.95"
$ws.Range("E12").Value = "  +2.59%  "

$ws.Range("E13").Value = "  +1.73%  "

$ws.Range("D14").Value = "0.889"
$ws.Range("E14").Value = "  +3.30%  "

$ws.Range("D15").Value = "2.628.38"
$ws.Range("E15").Value = "  +2.90%  "

$ws.Range("D16").Value = "14.70"
$ws.Range("E16").Value = "  +3.75%  "

$ws.Range("D17").Value = "2.281.28"
$ws.Range("E17").Value = "  +3.20%  "

$ws.Range("D18").Value = "44.111.86"
$ws.Range("E18").Value = "  +3.00%  "

$ws.Range("D19").Value = "14.11"
$ws.Range("E19").Value = "  -6.38%  "

$ws.Range("D20").Value = "0.0000100"
$ws.Range("E20").Value = "  +3.97%  "

$ws.Range("D21").Value = "6.58"
$ws.Range("E21").Value = "  +2.42%  "

$ws.Range("D22").Value = "66.55"
$ws.Range("E22").Value = "  +1.96%  "

$ws.Range("D23").Value = "3.22"
$ws.Range("E23").Value = "  +1.64%  "

$ws.Range("D24").Value = "238.69"
$ws.Range("E24").Value = "  +1.00%  "

$ws.Range("D25").Value = "2.21"
$ws.Range("E25").Value = "  +5.15%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").Value = "10.27"
$ws.Range("E27").Value = "  +1.99%  "

$ws.Range("D28").Value = "39.08"
$ws.Range("E28").Value = "  +14.82%  "

$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("D30").Value = "6.58"
$ws.Range("E30").Value = "  +4.16%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "20.70"
$ws.Range("E31").Value = "  +1.18%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "163.17"
$ws.Range("E32").Value = "  +4.05%  "

$ws.Range("D33").Value = "0.0891"
$ws.Range("E33").Value = "  +1.50%  "

$ws.Range("E34").Value = "  -1.21%  "

$ws.Range("E35").Value = "  +6.12%  "

$ws.Range("D36").Value = "3.28"
$ws.Range("E36").Value = "  +2.93%  "

$ws.Range("E37").Value = "  +12.27%  "

$ws.Range("D38").Value = "0.122"
$ws.Range("E38").Value = "  -0.19%  "

$ws.Range("E39").Value = "  +7.98%  "

$ws.Range("E40").Value = "  +1.06%  "

$ws.Range("D41").Value = "0.0330"
$ws.Range("E41").Value = "  +1.52%  "

$ws.Range("D42").Value = "15.47"
$ws.Range("E42").Value = "  +26.74%  "

$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("D44").Value = "1.765.94"
$ws.Range("E44").Value = "  -6.30%  "

$ws.Range("E45").Value = "  +1.13%  "

$ws.Range("D46").Value = "86.13"
$ws.Range("E46").Value = "  -3.30%  "

$ws.Range("D47").Value = "5.43"
$ws.Range("E47").Value = "  +0.49%  "

$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").Value = "75.62"
$ws.Range("E48").Value = "  +0.57%  "

$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "60.07"
$ws.Range("E49").Value = "  -0.75%  "

$ws.Range("D50").Value = "8.80"
$ws.Range("E50").Value = "  +2.71%  "

$ws.Range("D51").Value = "1.71"
$ws.Range("E51").Value = "  +6.40%  "
